$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value pairs derived from the authoritative diff
# (Price and Volume(1h) columns updated with refreshed market data)
$updates = @(
    @("D2", "278.75"),
    @("E2", "6.86%"),
    @("D3", "27.29"),
    @("E3", "0.85%"),
    @("D4", "4.800"),
    @("E4", "2.17%"),
    @("D5", "0.06249"),
    @("D6", "6.854"),
    @("E6", "1.82%"),
    @("D7", "3.264"),
    @("E7", "2.90%"),
    @("D8", "0.8769"),
    @("E8", "2.95%"),
    @("D9", "0.9511"),
    @("E9", "4.57%"),
    @("D10", "0.1444"),
    @("E10", "3.10%"),
    @("D11", "0.05089"),
    @("E11", "4.33%"),
    @("D12", "0.07277"),
    @("D13", "0.03158"),
    @("E13", "0.96%"),
    @("D14", "0.09041"),
    @("E14", "-0.18%"),
    @("D15", "0.001567"),
    @("E15", "2.18%"),
    @("D16", "0.0006267"),
    @("E16", "1.99%"),
    @("D17", "0.005918"),
    @("E17", "-2.76%"),
    @("E18", "0.51%"),
    @("E19", "5.09%"),
    @("D21", "0.1310"),
    @("E21", "0.09%"),
    @("D22", "3.865"),
    @("E22", "-5.81%"),
    @("D23", "0.04328"),
    @("E23", "1.77%"),
    @("E24", "-3.75%"),
    @("D25", "0.004274"),
    @("E25", "4.71%"),
    @("D26", "0.0001199"),
    @("E26", "-0.13%"),
    @("D27", "0.0001614"),
    @("E27", "-1.52%"),
    @("D40", "0.04044"),
    @("E40", "2.61%"),
    @("D41", "0.006713"),
    @("E41", "62.71%"),
    @("E42", "3.88%"),
    @("D43", "0.002208"),
    @("E43", "2.64%"),
    @("D44", "0.01411"),
    @("E44", "8.31%"),
    @("D45", "0.00005138"),
    @("E45", "0.28%"),
    @("E46", "-0.07%"),
    @("D47", "2.183"),
    @("E47", "3,228.21%"),
    @("E48", "-12.15%"),
    @("D49", "0.00002098"),
    @("E49", "-0.07%"),
    @("D50", "0.0001999"),
    @("E50", "-0.07%")
)

foreach ($pair in $updates) {
    $cellRef = $pair[0]
    $newValue = $pair[1]
    $rng = $ws.Range($cellRef)
    # Force text storage so numeric-looking strings and "%" literals
    # are preserved exactly as text, matching the original inline-string cells.
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
}
